# estimate_list_template.xlsx edit
#
# Commit intent: DB columns UNIT_TM_AMT / UNIT_HEAT_AMT were dropped and
# replaced by UNIT_MATERIAL_FINISH_TM_AMT / UNIT_MATERIAL_FINISH_HEAT_AMT,
# so the two merge-field placeholders that reference them need to be
# updated to match the new column names. A couple of small sheet-view /
# print-setup tweaks came along with the re-save as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K3 / K4: rename the merge-field placeholders -------------------------
# K3 used to render "${data.UNIT_TM_AMT}"   -> now "${data.UNIT_MATERIAL_FINISH_TM_AMT}"
# K4 used to render "${data.UNIT_HEAT_AMT}" -> now "${data.UNIT_MATERIAL_FINISH_HEAT_AMT}"
$ws.Range("K3").Value = '${data.UNIT_MATERIAL_FINISH_TM_AMT}'
$ws.Range("K4").Value = '${data.UNIT_MATERIAL_FINISH_HEAT_AMT}'

# --- Selection moved back to the top-left merged header cell --------------
[void]$ws.Range("A1:A2").Select()

# --- Print scale tightened from 50% to 46% ---------------------------------
$ws.PageSetup.Zoom = 46
